# Update "想去人数" (F column) figures on the 展览, 演出 and 全部类型 sheets
# to reflect the latest scrape output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 sheet -----------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 10519
$ws.Range("F6").Value  = 288
$ws.Range("F9").Value  = 786
$ws.Range("F13").Value = 3264
$ws.Range("F16").Value = 2187
$ws.Range("F17").Value = 2187
$ws.Range("F22").Value = 591
$ws.Range("F34").Value = 34
$ws.Range("F36").Value = 277
$ws.Range("F39").Value = 523
$ws.Range("F40").Value = 490
$ws.Range("F41").Value = 1739
$ws.Range("F44").Value = 59
$ws.Range("F45").Value = 479
$ws.Range("F46").Value = 1058

# --- 演出 sheet ------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 50

# --- 全部类型 sheet ---------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 10519
$ws.Range("F8").Value  = 288
$ws.Range("F11").Value = 786
$ws.Range("F13").Value = 3264
$ws.Range("F15").Value = 2187
$ws.Range("F16").Value = 2187
$ws.Range("F18").Value = 591
$ws.Range("F30").Value = 34
$ws.Range("F31").Value = 50
$ws.Range("F35").Value = 277
$ws.Range("F37").Value = 523
$ws.Range("F39").Value = 490
$ws.Range("F40").Value = 1739
$ws.Range("F46").Value = 59
$ws.Range("F47").Value = 479
$ws.Range("F48").Value = 1058
